$wb = $excel.ActiveWorkbook

# Rename sheets
$wsExio = $wb.Worksheets.Item("exiobase")
$wsExio.Name = "Exiobase"

$wsGerman = $wb.Worksheets.Item("german")
$wsGerman.Name = "Deutsch"

# Activate the "Deutsch" sheet and select cell C32 (linking settings tab to selection tab)
$wsGerman.Activate()
$wsGerman.Range("C32").Select()
